$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20, pushing existing rows 20:92 down to 21:93.
$ws.Rows("20:20").Insert()

# Populate the newly inserted row 20 with its data.
$ws.Range("A20").Value2 = 5
$ws.Range("B20").Value2 = "Macroferia Regional de Talca"
$ws.Range("C20").Value2 = "Maule"
$ws.Range("D20").Value2 = 45114
$ws.Range("E20").Value2 = 7
$ws.Range("F20").Value2 = "Fruta"
$ws.Range("G20").Value2 = 100107
$ws.Range("H20").Value2 = "Otros"
$ws.Range("I20").Value2 = 100107001
$ws.Range("J20").Value2 = "Caqui"
$ws.Range("K20").Value2 = "Mankaki"
$ws.Range("L20").Value2 = "Primera"
$ws.Range("M20").Value2 = 400
$ws.Range("N20").Value2 = 18000
$ws.Range("O20").Value2 = 18000
$ws.Range("P20").Value2 = 18000
$ws.Range("Q20").Value2 = "`$/caja 18 kilos granel"
$ws.Range("R20").Value2 = "Región del Maule"
$ws.Range("S20").Value2 = 1000
$ws.Range("T20").Value2 = 18

# The date column uses the same numeric date format as the rest of column D.
$ws.Range("D20").NumberFormat = $ws.Range("D21").NumberFormat
